$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new column at G (pushes old G -> H, old H -> I, etc.)
# ---------------------------------------------------------------------------
$ws.Columns("G:G").Insert()

# ---------------------------------------------------------------------------
# 2) Header row (row 10) new / moved headers
# ---------------------------------------------------------------------------
$ws.Range("G10").Value = "costo de Compra(realice una compra)"
$ws.Range("J10").Value = "dias/Costo de compra"
$ws.Range("K10").Value = "costo total del Inv Promedio"
$ws.Range("L10").Value = "Dias promedio Invenario"
$ws.Range("B10:L10").RowHeight = 63

# Copy the style used by the other headers (C10..I10, style s=4) onto the
# new header cells G10, J10, K10, L10
$ws.Range("C10").Copy() | Out-Null
$ws.Range("G10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("J10:L10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("G10").Value = "costo de Compra(realice una compra)"
$ws.Range("J10").Value = "dias/Costo de compra"
$ws.Range("K10").Value = "costo total del Inv Promedio"
$ws.Range("L10").Value = "Dias promedio Invenario"

# ---------------------------------------------------------------------------
# 3) Fill in the new "G" column (dias/costo de compra) constant values and
#    the new J (dias), K (costo total), L (rotacion) columns / formulas for
#    each data row 11-40
# ---------------------------------------------------------------------------
for ($r = 11; $r -le 32; $r++) {
    $ws.Cells.Item($r, 7).Value = 12
}
for ($r = 33; $r -le 40; $r++) {
    $ws.Cells.Item($r, 7).Value = 13
}
$ws.Range("G11:G40").Style = "Millares"

for ($r = 11; $r -le 40; $r++) {
    $j = 10 + ($r - 11)
    $ws.Cells.Item($r, 10).Value = $j
    $ws.Range("K$r").Formula = "=I$r*J$r"
    $ws.Range("L$r").Formula = "=I$r/K$r"
}
$ws.Range("K11:K40").Style = "Millares"

# Rows 11-15: M column = 1/L
for ($r = 11; $r -le 15; $r++) {
    $ws.Range("M$r").Formula = "=1/L$r"
}
$ws.Range("M11:M15").NumberFormat = "0.000000"

# L11:L15 highlighted yellow, 5 decimal number format
$ws.Range("L11:L15").NumberFormat = "0.00000"
$ws.Range("L11:L15").Interior.Color = 65535

# Row 15 extra: N15 average of M11:M15
$ws.Range("N15").Formula = "=SUM(M11:M15)/5"

# Row 16 extra: M16, N16
$ws.Range("M16").Formula = "=SUM(L11:L15)/5"
$ws.Range("N16").Formula = "=5/M16"
$ws.Range("L16").ClearFormats()

# ---------------------------------------------------------------------------
# 4) Row 41: totals row
# ---------------------------------------------------------------------------
$ws.Range("F41").Value = "costo de lo vendido"
$ws.Range("F41").Style = "Millares"
$ws.Range("G41").Formula = "=SUMPRODUCT(G11:G40,F11:F40)"
$ws.Range("G41").Style = "Millares"
$ws.Range("K41").Formula = "=I41*J41"
$ws.Range("K41").Style = "Millares"
$ws.Range("L41").Formula = "=SUM(L11:L40)"

# Row 42: empty but styled K42
$ws.Range("K42").Style = "Millares"

# ---------------------------------------------------------------------------
# 5) Row 43 (the old "promedio del periodo" summary row, now shifted down)
# ---------------------------------------------------------------------------
$ws.Range("K43").Formula = "=K41/30"
$ws.Range("K43").Style = "Millares"

# ---------------------------------------------------------------------------
# 6) New summary block rows 45-48 and row 51
# ---------------------------------------------------------------------------
$ws.Range("E45:F45").Merge()
$ws.Range("E45").Value = "Costo Inv Promedio"
$ws.Range("E45:F45").Style = "Millares"
$ws.Range("E45:F45").HorizontalAlignment = -4108   # xlCenter
$ws.Range("E45:F45").Interior.Color = 65535
$ws.Range("H45").Formula = "=I43*AVERAGE(G11:G40)"
$ws.Range("H45").Style = "Millares"
$ws.Range("I45").Style = "Millares"

$ws.Range("E46:F46").Merge()
$ws.Range("E46").Value = "Dias Promedio"
$ws.Range("E46:F46").Style = "Millares"
$ws.Range("E46:F46").HorizontalAlignment = -4108
$ws.Range("E46:F46").WrapText = $true
$ws.Range("E46:F46").Interior.Color = 65535
$ws.Range("H46").Formula = "=H45*30/G41"
$ws.Range("H46").Style = "Millares"

$ws.Range("C47").Value = "costo de ventas"
$ws.Range("C47").Style = "Millares"
$ws.Range("D47").Formula = "=F41*C44"
$ws.Range("D47").Style = "Millares"
$ws.Range("E47:F47").Merge()
$ws.Range("E47").Value = "Rotacion Inv"
$ws.Range("E47:F47").Style = "Millares"
$ws.Range("E47:F47").HorizontalAlignment = -4108
$ws.Range("E47:F47").WrapText = $true
$ws.Range("E47:F47").Interior.Color = 65535
$ws.Range("H47").Formula = "=30/H46"
$ws.Range("H47").Style = "Millares"

$ws.Range("C48:E48").Merge()
$ws.Range("C48").Value = "suma de todo el producto que se vendio durante el rango dias"
$ws.Range("C48:E48").Style = "Millares"
$ws.Range("C48:E48").HorizontalAlignment = -4108
$ws.Range("C48:E48").WrapText = $true
$ws.Range("F48:H48").Style = "Millares"
$ws.Range("F48:H48").WrapText = $true
$ws.Rows("48:48").RowHeight = 63

$ws.Range("C51").Value = "costo venta diario"
$ws.Range("C51").Style = "Millares"

# ---------------------------------------------------------------------------
# 7) Column widths
# ---------------------------------------------------------------------------
$ws.Columns("F:F").ColumnWidth = 18
$ws.Columns("G:G").ColumnWidth = 18.375
$ws.Columns("H:H").ColumnWidth = 16.625
$ws.Columns("I:I").ColumnWidth = 13.5
$ws.Columns("M:M").ColumnWidth = 11.375

# ---------------------------------------------------------------------------
# 8) View state: scroll / selection
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("C15"), $false)
$ws.Range("F41").Select() | Out-Null
